# Generate Report for Handoff
# File "b.md" has completed its handoff cycle (zh-cn + de-de), so update its
# status from "Handed back: in sync with en-US" to "Ready for handoff" and
# record the new handoff timestamps / handoff-file names for row 3 (b.md)
# on the Overview sheet and on each locale sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the "b.md" row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-23 04:07:26"

# --- zh-cn sheet: row 3 is the "b.md" row ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-23 04:07:18"

# --- de-de sheet: row 3 is the "b.md" row ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-23 04:07:26"
